$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:G1): lower-cased / renamed column titles ---
$ws.Range("A1").Value = "industry"
$ws.Range("B1").Value = "unit"
$ws.Range("C1").Value = "process"
$ws.Range("D1").Value = "carbon (kg CO2 eq)"
$ws.Range("E1").Value = "ced (MJ)"
$ws.Range("F1").Value = "climate change (kg CO2 eq)"
$ws.Range("G1").Value = "region"

# --- Data rows: columns D/E/F are permuted (D<-oldE, E<-oldF) and F gets a newly derived value ---
$ws.Range("D2").Value = 0.6090901733333334
$ws.Range("E2").Value = 8.7311964
$ws.Range("F2").Value = 0.00001698305
$ws.Range("D3").Value = 0.3927333866666667
$ws.Range("E3").Value = 6.1718623
$ws.Range("F3").Value = 0.000010950449
$ws.Range("D4").Value = 0.027492712
$ws.Range("E4").Value = 0.43205196
$ws.Range("F4").Value = 0.00000076656975
$ws.Range("D5").Value = 0.3927333866666667
$ws.Range("E5").Value = 6.1718623
$ws.Range("F5").Value = 0.000010950449
$ws.Range("D6").Value = 0.4350644066666667
$ws.Range("E6").Value = 6.2365689
$ws.Range("F6").Value = 0.00001213075
$ws.Range("D7").Value = 0.40242154
$ws.Range("E7").Value = 6.433261
$ws.Range("F7").Value = 0.00001122058
$ws.Range("D8").Value = 0.03219372333333334
$ws.Range("E8").Value = 0.51466088
$ws.Range("F8").Value = 0.00000089764642
$ws.Range("D9").Value = 0.40242154
$ws.Range("E9").Value = 6.433261
$ws.Range("F9").Value = 0.00001122058
$ws.Range("D10").Value = 0.4776557533333334
$ws.Range("E10").Value = 8.2683033
$ws.Range("F10").Value = 0.00001331831
$ws.Range("D11").Value = 0.4776557533333334
$ws.Range("E11").Value = 8.2683033
$ws.Range("F11").Value = 0.00001331831
$ws.Range("D12").Value = 0.4952027333333334
$ws.Range("E12").Value = 8.1555387
$ws.Range("F12").Value = 0.000013807566
$ws.Range("D13").Value = 0.6278351533333334
$ws.Range("E13").Value = 11.852454
$ws.Range("F13").Value = 0.00001750571
$ws.Range("D14").Value = 0.4856114666666667
$ws.Range("E14").Value = 7.896754
$ws.Range("F14").Value = 0.000013540136
$ws.Range("D15").Value = 0.3559258933333334
$ws.Range("E15").Value = 7.8827229
$ws.Range("F15").Value = 0.0000099241584
$ws.Range("D16").Value = 0.09547371333333333
$ws.Range("E16").Value = 4.2027104
$ws.Range("F16").Value = 0.0000026620604
$ws.Range("D17").Value = 0.09567875333333334
$ws.Range("E17").Value = 4.2088212
$ws.Range("F17").Value = 0.0000026677775
$ws.Range("D18").Value = 0.2165983733333333
$ws.Range("E18").Value = 6.8566748
$ws.Range("F18").Value = 0.0000060393373
$ws.Range("D19").Value = 0.1490731666666667
$ws.Range("E19").Value = 5.6571804
$ws.Range("F19").Value = 0.0000041565555
$ws.Range("D20").Value = 0.1655981266666667
$ws.Range("E20").Value = 5.3276359
$ws.Range("F20").Value = 0.0000046173152
$ws.Range("D21").Value = 0.13515254
$ws.Range("E21").Value = 5.5003998
$ws.Range("F21").Value = 0.0000037684115
$ws.Range("D22").Value = 0.1348409733333333
$ws.Range("E22").Value = 5.5655969
$ws.Range("F22").Value = 0.0000037597242
$ws.Range("D23").Value = 0.2439658066666667
$ws.Range("E23").Value = 5.0179742
$ws.Range("F23").Value = 0.000006802414
$ws.Range("D24").Value = 0.08371876666666667
$ws.Range("E24").Value = 3.9150982
$ws.Range("F24").Value = 0.0000023343013
$ws.Range("D25").Value = 0.20001732
$ws.Range("E25").Value = 5.6472019
$ws.Range("F25").Value = 0.0000055770136
$ws.Range("D26").Value = 0.1244235733333333
$ws.Range("E26").Value = 4.8492472
$ws.Range("F26").Value = 0.0000034692594
$ws.Range("D27").Value = 0.1721981
$ws.Range("E27").Value = 4.7940196
$ws.Range("F27").Value = 0.0000048013399
$ws.Range("D28").Value = 0.02200797933333333
$ws.Range("E28").Value = 0.49882554
$ws.Range("F28").Value = 0.00000061364085
$ws.Range("D29").Value = 0.9465277333333333
$ws.Range("E29").Value = 15.322566
$ws.Range("F29").Value = 0.000026391704
$ws.Range("D30").Value = 0.2952197933333333
$ws.Range("E30").Value = 5.3256997
$ws.Range("F30").Value = 0.0000082315111
$ws.Range("D31").Value = 0.4902247666666667
$ws.Range("E31").Value = 8.8459483
$ws.Range("F31").Value = 0.000013668767
$ws.Range("D32").Value = 0.1326324133333333
$ws.Range("E32").Value = 3.6969153
$ws.Range("F32").Value = 0.0000036981437
$ws.Range("D33").Value = 0.01223404533333333
$ws.Range("E33").Value = 0.28401835
$ws.Range("F33").Value = 0.00000034111764
$ws.Range("D34").Value = 0.02200797933333333
$ws.Range("E34").Value = 0.49882554
$ws.Range("F34").Value = 0.00000061364085
$ws.Range("D35").Value = 0.007259669333333333
$ws.Range("E35").Value = 0.1310752
$ws.Range("F35").Value = 0.00000020241885
$ws.Range("D36").Value = 0.00943757
$ws.Range("E36").Value = 0.17039776
$ws.Range("F36").Value = 0.0000002631445
$ws.Range("D37").Value = 0.007948737333333334
$ws.Range("E37").Value = 0.20521385
$ws.Range("F37").Value = 0.00000022163188

# --- Header-cell comments describing each column's data type ---
$ws.Range("A1").AddComment("Data type: Categorical (text)") | Out-Null
$ws.Range("B1").AddComment("Data type: Various (e.g. kg, kWh)") | Out-Null
$ws.Range("C1").AddComment("Data type: Categorical (text)") | Out-Null
$ws.Range("D1").AddComment("Data type: Carbon footprint") | Out-Null
$ws.Range("E1").AddComment("Data type: Cumulative energy demand") | Out-Null
$ws.Range("F1").AddComment("Data type: Climate change impact") | Out-Null
$ws.Range("G1").AddComment("Data type: Categorical (text)") | Out-Null
